# "added columns to timesheet"
# Sheet1 gets a header row (date / start / stop / interrupt /
# delta (time-interrupt) / activity / comments) and the two wide
# text columns (E, G) are auto-fit to their header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "start"
$ws.Range("C1").Value = "stop"
$ws.Range("D1").Value = "interrupt"
$ws.Range("E1").Value = "delta (time-interrupt)"
$ws.Range("F1").Value = "activity"
$ws.Range("G1").Value = "comments"

# Widen the two columns whose header text doesn't fit the default width.
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null

# Leave the cursor parked under the header row, like a user who just
# finished typing the headers and tabbed/entered down to A2.
$ws.Range("A2").Select() | Out-Null
